# "error solve ifrs list"
#
# The per-year IFRS metrics for KC코트렐 (company_list sheet) were loaded at
# the wrong scale / from the wrong source columns. This corrects the actual
# figures for FY2014-FY2018 (rows 2-6) and blanks out the bogus forecast
# rows FY2019(E)-FY2021(E) (rows 7-9), which only keeps their index/label
# columns (A:C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected values per row/column for the FY2014-FY2018 (rows 2-6) data.
# Columns intentionally absent from a row here (e.g. J/O on row 2, or J/O on
# row 6) are removed below via ClearContents rather than given a value.
$rowData = @{
    2 = @{
        D = 2097;  E = -63;   F = -63;   G = -44;   H = -76;   I = -76
        K = 1707;  L = 1104;  M = 603;   N = 603
        P = 50;    Q = 46;    R = -60;   S = 10;    T = 60;    U = -14;   V = 135
        W = -3.03; X = -3.65; Y = -11.58; Z = -4.58; AA = 182.99; AB = 1110.07
        AC = -706; AD = -8.98; AE = 5566; AF = 1.14; AG = 0; AH = 0; AI = 0
        AJ = 10837872
    }
    3 = @{
        D = 1868;  E = -87;   F = -87;   G = -84;   H = -68;   I = -66;   J = -1
        K = 1467;  L = 907;   M = 561;   N = 555;   O = 6
        P = 50;    Q = -19;   R = -13;   S = -36;   T = 15;    U = -34;   V = 100
        W = -4.67; X = -3.61; Y = -11.42; Z = -4.25; AA = 161.69; AB = 1007.37
        AC = -610; AD = -10.77; AE = 5121; AF = 1.28; AG = 0; AH = 0; AI = 0
        AJ = 10837872
    }
    4 = @{
        D = 1793;  E = 60;    F = 60;    G = 12;    H = 9;     I = 9;     J = 0
        K = 1535;  L = 832;   M = 703;   N = 690;   O = 13
        P = 65;    Q = 119;   R = -66;   S = 63;    T = 5;     U = 114;   V = 26
        W = 3.33;  X = 0.51;  Y = 1.43;  Z = 0.61;  AA = 118.3; AB = 959.38
        AC = 75;   AD = 79.45; AE = 5305; AF = 1.12; AG = 40; AH = 0.67; AI = 58.43
        AJ = 13000000
    }
    5 = @{
        D = 1685;  E = 35;    F = 35;    G = 93;    H = 69;    I = 67;    J = 2
        K = 1803;  L = 1026;  M = 777;   N = 762;   O = 15
        P = 65;    Q = 30;    R = -16;   S = -9;    T = 20;    U = 10;    V = 45
        W = 2.05;  X = 4.09;  Y = 9.26;  Z = 4.13;  AA = 132;  AB = 1071.92
        AC = 517;  AD = 9.26; AE = 5861; AF = 0.82; AG = 100; AH = 2.09; AI = 19.33
        AJ = 13000000
    }
    6 = @{
        D = 1820;  E = 37;    F = 37;    G = 27;    H = 22;    I = 42
        K = 2195;  L = 1417;  M = 777;   N = 780
        P = 65;    Q = -68;   R = -63;   S = 113;   T = 5;     U = -73;   V = 136
        W = 2.05;  X = 1.21;  Y = 5.41;  Z = 1.1;   AA = 182.36; AB = 1088.54
        AC = 320;  AD = 13.29; AE = 5996; AF = 0.71; AG = 80; AH = 1.88; AI = 24.97
        AJ = 13000000
    }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}

# J2 and O2 must be removed entirely (not merely set to 0/blank).
foreach ($addr in @("J2", "O2")) {
    $ws.Range($addr).ClearContents()
}

# Rows 7-9 (2019(E)/2020(E)/2021(E) forecast rows) were erroneous in full -
# wipe every data column (D:AJ), keeping only the index/annual/period
# columns A:C.
foreach ($r in 7..9) {
    $ws.Range("D${r}:AJ${r}").ClearContents()
}
